# Aggregate to monthly working
# - Fix "Acre-feet" -> "acre-feet" on the Fields sheet (Storage + Evaporation rows)
# - Add a new "AggregateByTimePeriod" column (D) to the Fields sheet with Yes/No values
# - Make "Fields" the active sheet/tab with D1 selected

$wb = $excel.ActiveWorkbook

$wsFields = $wb.Worksheets.Item("Fields")
$wsReservoirs = $wb.Worksheets.Item("Reservoirs")

# Normalize the units text for Storage / Evaporation to lowercase "acre-feet"
$wsFields.Range("C2").Value = "acre-feet"
$wsFields.Range("C3").Value = "acre-feet"

# New header
$wsFields.Range("D1").Value = "AggregateByTimePeriod"

# New column values per FieldID row
$wsFields.Range("D2").Value = "No"   # 17 Storage
$wsFields.Range("D3").Value = "Yes"  # 25 Evaporation
$wsFields.Range("D4").Value = "Yes"  # 29 Inflow
$wsFields.Range("D5").Value = "Yes"  # 30 Inflow Volume
$wsFields.Range("D6").Value = "Yes"  # 33 Unregulated Inflow
$wsFields.Range("D7").Value = "Yes"  # 34 Unregulated Inflow Volume
$wsFields.Range("D8").Value = "Yes"  # 42 Total Release
$wsFields.Range("D9").Value = "Yes"  # 43 Release volume
$wsFields.Range("D10").Value = "No"  # 49 Pool Elevation
$wsFields.Range("D11").Value = "No"  # 89 Area

# Make Fields the active/selected sheet, with D1 selected
$wsFields.Activate()
$wsFields.Range("D1").Select()
